# Generate Report for Archive
# Update the "Status" value from "Ready for handoff" to "In Translation"
# across the Overview sheet and each language sheet, then let the
# column widths re-flow (autofit) to match the shorter text.

$wb = $excel.ActiveWorkbook

# Target column width, matching what Excel's AutoFit settles on once the
# "Status" columns hold the shorter "In Translation" text instead of
# "Ready for handoff". (ColumnWidth is quantized by the host to 1/6-wide
# steps, so 12.5 is the closest settable value that lands on this target.)
$targetStatusColWidth = 12.5

# --- Overview sheet: zh-cn / de-de status columns (E & F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRange = $wsOverview.Range("E2:F4")
for ($r = 1; $r -le $overviewRange.Rows.Count; $r++) {
    for ($c = 1; $c -le $overviewRange.Columns.Count; $c++) {
        $cell = $overviewRange.Cells.Item($r, $c)
        if ($cell.Text -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}
$wsOverview.Range("E:F").EntireColumn.AutoFit() | Out-Null
$wsOverview.Range("E1").EntireColumn.ColumnWidth = $targetStatusColWidth
$wsOverview.Range("F1").EntireColumn.ColumnWidth = $targetStatusColWidth

# --- Language sheets: Status column (C) ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $statusRange = $ws.Range("C2:C4")
    for ($r = 1; $r -le $statusRange.Rows.Count; $r++) {
        $cell = $statusRange.Cells.Item($r, 1)
        if ($cell.Text -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
    $ws.Range("C:C").EntireColumn.AutoFit() | Out-Null
    $ws.Range("C1").EntireColumn.ColumnWidth = $targetStatusColWidth
}
